$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured value in B3; the dependent AVERAGE in B32 recalculates automatically.
$ws.Range("B3").Value = 0.9299

# Scroll the view back so the top-left visible cell is A1 (clears the saved scroll position).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
